# Update the "想去人数" (F column) numbers on the "展览" and "全部类型" sheets
# to reflect the newly generated data (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 233
$ws1.Range("F4").Value = 4829
$ws1.Range("F5").Value = 214
$ws1.Range("F6").Value = 0
$ws1.Range("F8").Value = 108
$ws1.Range("F14").Value = 258
$ws1.Range("F18").Value = 153
$ws1.Range("F20").Value = 4046
$ws1.Range("F21").Value = 6367
$ws1.Range("F24").Value = 88
$ws1.Range("F31").Value = 2589
$ws1.Range("F33").Value = 534
$ws1.Range("F35").Value = 298
$ws1.Range("F38").Value = 179
$ws1.Range("F42").Value = 47
$ws1.Range("F46").Value = 0
$ws1.Range("F49").Value = 587

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 233
$ws4.Range("F4").Value = 4829
$ws4.Range("F8").Value = 108
$ws4.Range("F13").Value = 115
$ws4.Range("F14").Value = 258
$ws4.Range("F18").Value = 153
$ws4.Range("F19").Value = 0
$ws4.Range("F20").Value = 0
$ws4.Range("F21").Value = 6367
$ws4.Range("F27").Value = 0
$ws4.Range("F31").Value = 2589
$ws4.Range("F35").Value = 298
$ws4.Range("F37").Value = 376
$ws4.Range("F38").Value = 179
$ws4.Range("F39").Value = 0
$ws4.Range("F40").Value = 0
$ws4.Range("F45").Value = 497
$ws4.Range("F48").Value = 78
